$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the existing "Yht" summary row (currently row 17) down to row 19,
#    keeping its formatting, label and formula (updated range) intact.
# ---------------------------------------------------------------------------
$ws.Range("B17:D17").Copy()
$ws.Range("B19:D19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B19").Value = "Yht"
$ws.Range("C19").Formula = "=SUM(C6:C17)"
$ws.Range("D19").Value = ""
$ws.Rows("19:19").RowHeight = 18.75

# ---------------------------------------------------------------------------
# 2) Turn the old row 17 into a regular data row (same look as row 16) and
#    add a new blank spacer row 18 under it, using row 16 as the format
#    template so the same style indices get reused.
# ---------------------------------------------------------------------------
$ws.Range("B16:D16").Copy()
$ws.Range("B17:D18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Clear any stale values that used to live in row 17 (old "Yht" row) / row 18
$ws.Range("B17:D18").ClearContents()

# ---------------------------------------------------------------------------
# 3) Fill in the new diary entry on row 17.
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = 45335
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = "Koko tämän ajan käytin yritääkseni saamaa lisätyä uusia täppiä navigaatio palkkiin erillään muista. Tämä johti lukemaan dokumentaatiota ja ihmetelemään tekoälyn tyhmiä vastauksia. Lopulta sain todella hyvän lopputuloksen mutta tämä on jo toinen kerta kun navigaatio palkki aiheutti näin paljon pään särkyä saman projektin aikana. Myöskin kirjoitin about sivun."

$ws.Rows("17:17").RowHeight = 131.25
$ws.Rows("18:18").RowHeight = 18.75

# ---------------------------------------------------------------------------
# 4) Update the view: scroll so row 13 is at the top and select D17, matching
#    where the author was working when they saved.
# ---------------------------------------------------------------------------
$ws.Range("D17").Select()
